# Update the FHIR ValueSet metadata workbook to point at the new
# LinuxForHealth URLs/publisher/version/date instead of the old
# ibm.com/Alvearie values.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/match-method"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Include-from sheet ---
$inc = $wb.Worksheets.Item("Include from Match Method Cod")

# System URI
$inc.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/match-method"
